# =====================================================================
# feat: add 2022-Q3 data
#
# - Insert a new "2022-Q3" worksheet (cloned from "2022-Q2" so it keeps the
#   same header/border styling) right after the totals sheet, and fill it
#   with the Q3 per-fund holdings.
# - Insert a new row at the top of the "..." (totals) sheet data block
#   summarising the new quarter, renumbering the trailing index column.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: add & populate the "2022-Q3" worksheet
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)            # existing "2022-Q2" sheet (template)
$q2Sheet.Copy($q2Sheet)                       # duplicate; copy is placed right before it
$q3Sheet = $wb.Worksheets.Item(2)             # the freshly inserted copy
$q3Sheet.Name = "2022-Q3"

# The template has 10 data rows (rows 2-11); this quarter needs 14
# data rows (rows 2-15). Insert the 4 extra rows, then paste the last
# template rows formatting onto them so new rows keep identical styling.
$q3Sheet.Rows.Item(11).Insert()
$q3Sheet.Rows.Item(11).Insert()
$q3Sheet.Rows.Item(11).Insert()
$q3Sheet.Rows.Item(11).Insert()
$q3Sheet.Range("A10:H10").Copy()
$q3Sheet.Range("A11:H14").PasteSpecial(-4122)
$q3Sheet.Application.CutCopyMode = $false

# Columns B, D, E, F, G hold numeric-looking values that must stay text
# (leading zeros in fund codes, fixed trailing-zero decimals, etc.) - format
# the whole data block as Text up front so every row shares one style.
$q3Sheet.Range("B2:B15").NumberFormat = "@"
$q3Sheet.Range("D2:G15").NumberFormat = "@"

# Header row
$q3Sheet.Cells.Item(1,2).Value = "基金代码"
$q3Sheet.Cells.Item(1,3).Value = "基金名称"
$q3Sheet.Cells.Item(1,4).Value = "基金规模"
$q3Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q3Sheet.Cells.Item(1,6).Value = "仓位占比"
$q3Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3Sheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows: A = index, B = fund code, C = fund name, D = fund size,
# E = total stock position, F = position share, G = held market value,
# H = position rank
$q3Sheet.Cells.Item(2,1).Value = 0
$q3Sheet.Cells.Item(2,2).Value = "000043"
$q3Sheet.Cells.Item(2,3).Value = "嘉实美国成长股票（QDII）人民币"
$q3Sheet.Cells.Item(2,4).Value = "12.41"
$q3Sheet.Cells.Item(2,5).Value = "92.80"
$q3Sheet.Cells.Item(2,6).Value = "1.31"
$q3Sheet.Cells.Item(2,7).Value = "0.1626"
$q3Sheet.Cells.Item(2,8).Value = 9

$q3Sheet.Cells.Item(3,1).Value = 1
$q3Sheet.Cells.Item(3,2).Value = "000044"
$q3Sheet.Cells.Item(3,3).Value = "嘉实美国成长股票（QDII）美元现汇"
$q3Sheet.Cells.Item(3,4).Value = "12.41"
$q3Sheet.Cells.Item(3,5).Value = "92.80"
$q3Sheet.Cells.Item(3,6).Value = "1.31"
$q3Sheet.Cells.Item(3,7).Value = "0.1626"
$q3Sheet.Cells.Item(3,8).Value = 9

$q3Sheet.Cells.Item(4,1).Value = 2
$q3Sheet.Cells.Item(4,2).Value = "000369"
$q3Sheet.Cells.Item(4,3).Value = "广发全球医疗保健（QDII）人民币A"
$q3Sheet.Cells.Item(4,4).Value = "2.76"
$q3Sheet.Cells.Item(4,5).Value = "83.19"
$q3Sheet.Cells.Item(4,6).Value = "3.85"
$q3Sheet.Cells.Item(4,7).Value = "0.1063"
$q3Sheet.Cells.Item(4,8).Value = 3

$q3Sheet.Cells.Item(5,1).Value = 3
$q3Sheet.Cells.Item(5,2).Value = "000370"
$q3Sheet.Cells.Item(5,3).Value = "广发全球医疗保健（QDII）美元A"
$q3Sheet.Cells.Item(5,4).Value = "2.75"
$q3Sheet.Cells.Item(5,5).Value = "83.19"
$q3Sheet.Cells.Item(5,6).Value = "3.85"
$q3Sheet.Cells.Item(5,7).Value = "0.1059"
$q3Sheet.Cells.Item(5,8).Value = 3

$q3Sheet.Cells.Item(6,1).Value = 4
$q3Sheet.Cells.Item(6,2).Value = "161126"
$q3Sheet.Cells.Item(6,3).Value = "易方达标普医疗保健指数（QDII-LOF）人民币"
$q3Sheet.Cells.Item(6,4).Value = "0.57"
$q3Sheet.Cells.Item(6,5).Value = "93.29"
$q3Sheet.Cells.Item(6,6).Value = "1.65"
$q3Sheet.Cells.Item(6,7).Value = "0.0094"
$q3Sheet.Cells.Item(6,8).Value = 2

$q3Sheet.Cells.Item(7,1).Value = 5
$q3Sheet.Cells.Item(7,2).Value = "012864"
$q3Sheet.Cells.Item(7,3).Value = "易方达标普医疗保健指数（QDII-LOF）人民币 C"
$q3Sheet.Cells.Item(7,4).Value = "0.57"
$q3Sheet.Cells.Item(7,5).Value = "93.29"
$q3Sheet.Cells.Item(7,6).Value = "1.65"
$q3Sheet.Cells.Item(7,7).Value = "0.0094"
$q3Sheet.Cells.Item(7,8).Value = 2

$q3Sheet.Cells.Item(8,1).Value = 6
$q3Sheet.Cells.Item(8,2).Value = "003719"
$q3Sheet.Cells.Item(8,3).Value = "易方达标普医疗保健指数（QDII-LOF）美元A"
$q3Sheet.Cells.Item(8,4).Value = "0.45"
$q3Sheet.Cells.Item(8,5).Value = "93.29"
$q3Sheet.Cells.Item(8,6).Value = "1.65"
$q3Sheet.Cells.Item(8,7).Value = "0.0074"
$q3Sheet.Cells.Item(8,8).Value = 2

$q3Sheet.Cells.Item(9,1).Value = 7
$q3Sheet.Cells.Item(9,2).Value = "096001"
$q3Sheet.Cells.Item(9,3).Value = "大成标普500等权重指数（QDII）人民币"
$q3Sheet.Cells.Item(9,4).Value = "3.08"
$q3Sheet.Cells.Item(9,5).Value = "93.16"
$q3Sheet.Cells.Item(9,6).Value = "0.22"
$q3Sheet.Cells.Item(9,7).Value = "0.0068"
$q3Sheet.Cells.Item(9,8).Value = 5

$q3Sheet.Cells.Item(10,1).Value = 8
$q3Sheet.Cells.Item(10,2).Value = "013404"
$q3Sheet.Cells.Item(10,3).Value = "大成标普500等权重指数（QDII）美元"
$q3Sheet.Cells.Item(10,4).Value = "3.08"
$q3Sheet.Cells.Item(10,5).Value = "93.16"
$q3Sheet.Cells.Item(10,6).Value = "0.22"
$q3Sheet.Cells.Item(10,7).Value = "0.0068"
$q3Sheet.Cells.Item(10,8).Value = 5

$q3Sheet.Cells.Item(11,1).Value = 9
$q3Sheet.Cells.Item(11,2).Value = "011706"
$q3Sheet.Cells.Item(11,3).Value = "长信美国标准普尔100等权重指数增强（QDII）美元"
$q3Sheet.Cells.Item(11,4).Value = "0.39"
$q3Sheet.Cells.Item(11,5).Value = "82.64"
$q3Sheet.Cells.Item(11,6).Value = "0.92"
$q3Sheet.Cells.Item(11,7).Value = "0.0036"
$q3Sheet.Cells.Item(11,8).Value = 1

$q3Sheet.Cells.Item(12,1).Value = 10
$q3Sheet.Cells.Item(12,2).Value = "519981"
$q3Sheet.Cells.Item(12,3).Value = "长信美国标准普尔100等权重指数增强（QDII）人民币"
$q3Sheet.Cells.Item(12,4).Value = "0.39"
$q3Sheet.Cells.Item(12,5).Value = "82.64"
$q3Sheet.Cells.Item(12,6).Value = "0.92"
$q3Sheet.Cells.Item(12,7).Value = "0.0036"
$q3Sheet.Cells.Item(12,8).Value = 1

$q3Sheet.Cells.Item(13,1).Value = 11
$q3Sheet.Cells.Item(13,2).Value = "012865"
$q3Sheet.Cells.Item(13,3).Value = "易方达标普医疗保健指数（QDII-LOF）美元 C"
$q3Sheet.Cells.Item(13,4).Value = "0.12"
$q3Sheet.Cells.Item(13,5).Value = "93.29"
$q3Sheet.Cells.Item(13,6).Value = "1.65"
$q3Sheet.Cells.Item(13,7).Value = "0.0020"
$q3Sheet.Cells.Item(13,8).Value = 2

$q3Sheet.Cells.Item(14,1).Value = 12
$q3Sheet.Cells.Item(14,2).Value = "016280"
$q3Sheet.Cells.Item(14,3).Value = "广发全球医疗保健（QDII）人民币C"
$q3Sheet.Cells.Item(14,4).Value = "0.02"
$q3Sheet.Cells.Item(14,5).Value = "83.19"
$q3Sheet.Cells.Item(14,6).Value = "3.85"
$q3Sheet.Cells.Item(14,7).Value = "0.0008"
$q3Sheet.Cells.Item(14,8).Value = 3

$q3Sheet.Cells.Item(15,1).Value = 13
$q3Sheet.Cells.Item(15,2).Value = "016281"
$q3Sheet.Cells.Item(15,3).Value = "广发全球医疗保健（QDII）美元C"
$q3Sheet.Cells.Item(15,4).Value = "0.02"
$q3Sheet.Cells.Item(15,5).Value = "83.19"
$q3Sheet.Cells.Item(15,6).Value = "3.85"
$q3Sheet.Cells.Item(15,7).Value = "0.0008"
$q3Sheet.Cells.Item(15,8).Value = 3

# ---------------------------------------------------------------------
# Part 2: update the totals sheet with the new quarters summary row
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Rows.Item(2).Insert()                  # push existing data rows down by one
$totalSheet.Range("A3:D3").Copy()                  # copy formatting of a data row onto the
$totalSheet.Range("A2:D2").PasteSpecial(-4122)      # blank new row 2 (keeps column-A style)
$totalSheet.Application.CutCopyMode = $false

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q3"
$totalSheet.Cells.Item(2,3).Value = 14
$totalSheet.Cells.Item(2,4).Value = 0.59

# Column A is a plain sequential index (row number - 2); refresh it for every
# row now that a new row has been inserted at the top.
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(7,1).Value = 5

